# SMP_4PP.xlsx edit script
# Relates characteristic state line with deviatoric strain:
# inserts two new columns (AQ, AR) holding "dratio-M" and "dratio_pred",
# shifting the old AQ:AV block to AS:AX, and fixes up the formulas that
# referenced the old AQ column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new columns before the old AQ column (column 43).
#    Excel automatically shifts the existing AQ:AV content (and any
#    formula references to it) two columns to the right, to AS:AX.
$ws.Range("AQ1:AR1").EntireColumn.Insert()

# 2) Header labels for the two new columns.
$ws.Range("AQ1").Value = "dratio-M"
$ws.Range("AR1").Value = "dratio_pred"

# 3) Fill the new AQ / AR columns (rows 2-18) with the new formulas, and
#    repair the AN column formula so that it refers to the (now shifted)
#    AS column instead of the old AQ column.
for ($r = 2; $r -le 18; $r++) {
    $ws.Range("AQ$r").Formula = "=AL$r/AK$r-(1.35*(AK$r/3255000)^-0.0723)"
    $ws.Range("AR$r").Formula = "=1.35*(AK$r/3255000)^-0.0723-1.386/(AE$r+1.27)+0.03463"
    $ws.Range("AN$r").Formula = "=1/(2+AM$r*AS$r-2*0.33*(1+AM$r+AS$r))"
}

# 4) Column widths: the two freshly inserted columns get a width of
#    their own, distinct from the columns that used to sit at AQ/AR.
$ws.Range("AQ1:AR1").EntireColumn.ColumnWidth = 9.65

# 5) Sheet view bookkeeping - the author had scrolled / selected a
#    different column after making the edit.
$ws.Activate()
$ws.Range("AS1:AS1048576").Select()
$excel.ActiveWindow.ScrollColumn = 32
$excel.ActiveWindow.ScrollRow = 1
